# organisations: Add reference in project data, import updates to
# organisation records, fill data when viewing projects.
#
# Adds a new "Outcome Fund" worksheet (after "Outcomes") that mirrors the
# existing SPREADSHEETFORM:DOWN layout used by "Outcomes", and nudges the
# remembered cell-selection on the first two sheets.

$wb = $excel.ActiveWorkbook

$wsGeneral  = $wb.Worksheets.Item("General Overview")
$wsOutcomes = $wb.Worksheets.Item("Outcomes")

# ---------------------------------------------------------------------
# 1. New "Outcome Fund" sheet, inserted right after "Outcomes".
# ---------------------------------------------------------------------
$wsFund = $wb.Worksheets.Add($null, $wsOutcomes)
$wsFund.Name = "Outcome Fund"

# Column widths (character units = LibreOffice width(chars) - 5/6, which is
# the padding Excel adds on top of the stored character width).
$wsFund.Columns.Item(1).ColumnWidth = 23.6167
$wsFund.Columns.Item(2).ColumnWidth = 22.5167
$wsFund.Columns.Item(3).ColumnWidth = 21.7967
$wsFund.Columns.Item(4).ColumnWidth = 21.2567
$wsFund.Columns.Item(5).ColumnWidth = 19.9967

# Merge the banner cells *before* applying formatting -- merging after
# formatting splits the thin border into a distinct per-edge style, which
# would introduce extra cellXfs entries that the original file doesn't have.
$wsFund.Range("A3:B3").Merge() | Out-Null
$wsFund.Range("C3:E3").Merge() | Out-Null

# Reuse the "section title" style (fontId 6 / blue fill / thin border,
# general alignment) for row 3 (the "Fund" / "Organisation" banner) and
# row 4 (field labels), sourced from Outcomes!B4:F4 which already carries
# that exact style.
$wsOutcomes.Range("B4:F4").Copy() | Out-Null
$wsFund.Range("A3:E3").PasteSpecial(-4122) | Out-Null
$wsFund.Range("A4:E4").PasteSpecial(-4122) | Out-Null

# Reuse the plain data-row style (grey fill, thin border) for rows 5-17,
# sourced from Outcomes!A6:E6 (an already-blank row in that style).
$wsOutcomes.Range("A6:E6").Copy() | Out-Null
for ($r = 5; $r -le 17; $r++) {
    $wsFund.Range("A" + $r + ":E" + $r).PasteSpecial(-4122) | Out-Null
}

# A4 ("Name") is left-aligned rather than general -- this introduces the
# one genuinely new cell style.
$wsFund.Range("A4").HorizontalAlignment = -4131

# Header banner + field labels + SPREADSHEETFORM source row.
$wsFund.Range("A3").Value = "Fund"
$wsFund.Range("C3").Value = "Organisation"

$wsFund.Range("A4").Value = "Name"
$wsFund.Range("B4").Value = "Identifier"
$wsFund.Range("C4").Value = "Our ID"
$wsFund.Range("D4").Value = "Name"
$wsFund.Range("E4").Value = "Type"

$wsFund.Range("A5").Value = "SPREADSHEETFORM:DOWN:outcome_funds:title"
$wsFund.Range("B5").Value = "SPREADSHEETFORM:DOWN:outcome_funds:id"
$wsFund.Range("C5").Value = "SPREADSHEETFORM:DOWN:outcome_funds:organisation/id"
$wsFund.Range("D5").Value = "SPREADSHEETFORM:DOWN:outcome_funds:organisation/name"
$wsFund.Range("E5").Value = "SPREADSHEETFORM:DOWN:outcome_funds:organisation/type"

$wsFund.Range("A2").Select() | Out-Null

# ---------------------------------------------------------------------
# 2. Remembered selection tweaks on the existing sheets.
# ---------------------------------------------------------------------
$wsGeneral.Range("A23").Select() | Out-Null
$wsOutcomes.Range("A1").Select() | Out-Null

# Restore "General Overview" as the active tab (adding/selecting on the
# other sheets shifts focus as a side effect).
$wsGeneral.Activate() | Out-Null
